# Horas.xlsx - "Estudio de alternativas. Anteproyecto"
# Append three new work-log entries (rows 20-22) for Federico Speroni,
# matching the rows already present (Nombre, Fecha, Horas, Nombre de la
# tarea realizada, Detalle).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20: 24/04/2017 - 4 hs - Creación de Anteproyecto
$ws.Range("A20").Value = "Federico Speroni"
$ws.Range("B20").Value = 42849
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = "Creación de Anteproyecto"
$ws.Range("E20").Value = "Retoques Anteproyecto, gráfico de Gantt"

# Row 21: 25/04/2017 - 4 hs - Creación de Anteproyecto
$ws.Range("A21").Value = "Federico Speroni"
$ws.Range("B21").Value = 42850
$ws.Range("C21").Value = 4
$ws.Range("D21").Value = "Creación de Anteproyecto"
$ws.Range("E21").Value = "Retoques Anteproyecto, Estudio de alternativas"

# Row 22: 26/04/2017 - 1 hs - Investigación
$ws.Range("A22").Value = "Federico Speroni"
$ws.Range("B22").Value = 42851
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = "Investigación"
$ws.Range("E22").Value = "Video tutorial Angular 2"

# Match the date formatting already used for column B (numFmtId 14)
# by copying the format from the row immediately above the new ones.
$ws.Range("B19").Copy()
$ws.Range("B20:B22").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Leave the selection where the author left it after typing the data.
$ws.Range("C22").Select()
